# Excess mortality analyses - Week 44
# Updates several provinces' weekly death counts (RIVM refresh) and adds the
# newly published weeks 43 & 44 (rows 149 & 150) with their variance formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 119 ---------------------------------------------------------------
$ws.Range("X119").Value = 724

# --- Row 126 ---------------------------------------------------------------
$ws.Range("X126").Value = 574

# --- Row 127 ---------------------------------------------------------------
$ws.Range("W127").Value = 434

# --- Row 132 ---------------------------------------------------------------
$ws.Range("U132").Value = 374

# --- Row 133 ---------------------------------------------------------------
$ws.Range("AA133").Value = 209

# --- Row 135 ---------------------------------------------------------------
$ws.Range("S135").Value = 207
$ws.Range("X135").Value = 640

# --- Row 136 ---------------------------------------------------------------
$ws.Range("X136").Value = 602

# --- Row 138 ---------------------------------------------------------------
$ws.Range("S138").Value = 212
$ws.Range("U138").Value = 395
$ws.Range("W138").Value = 448
$ws.Range("X138").Value = 598

# --- Row 139 ---------------------------------------------------------------
$ws.Range("V139").Value = 211
$ws.Range("X139").Value = 608
$ws.Range("AA139").Value = 206

# --- Row 140 ---------------------------------------------------------------
$ws.Range("X140").Value = 597

# --- Row 141 ---------------------------------------------------------------
$ws.Range("U141").Value = 396

# --- Row 142 ---------------------------------------------------------------
$ws.Range("X142").Value = 594

# --- Row 143 ---------------------------------------------------------------
$ws.Range("P143").Value = 111
$ws.Range("V143").Value = 170
$ws.Range("W143").Value = 406
$ws.Range("X143").Value = 605
$ws.Range("AA143").Value = 225

# --- Row 144 ---------------------------------------------------------------
$ws.Range("W144").Value = 487
$ws.Range("Z144").Value = 427

# --- Row 145 ---------------------------------------------------------------
$ws.Range("W145").Value = 485
$ws.Range("X145").Value = 660
$ws.Range("Y145").Value = 77

# --- Row 146 ---------------------------------------------------------------
$ws.Range("S146").Value = 246
$ws.Range("V146").Value = 217
$ws.Range("W146").Value = 500
$ws.Range("X146").Value = 662
$ws.Range("AA146").Value = 223

# --- Row 147 ---------------------------------------------------------------
$ws.Range("S147").Value = 232
$ws.Range("U147").Value = 459
$ws.Range("V147").Value = 215
$ws.Range("W147").Value = 475
$ws.Range("X147").Value = 689
$ws.Range("Z147").Value = 484
$ws.Range("AA147").Value = 257

# --- Row 148 ---------------------------------------------------------------
$ws.Range("P148").Value = 119
$ws.Range("Q148").Value = 150
$ws.Range("S148").Value = 251
$ws.Range("T148").Value = 63
$ws.Range("U148").Value = 419
$ws.Range("V148").Value = 237
$ws.Range("W148").Value = 496
$ws.Range("X148").Value = 657
$ws.Range("Z148").Value = 534
$ws.Range("AA148").Value = 227

# --- Row 149 (new week 2022-43) --------------------------------------------
$ws.Range("N149").Value = 2022
$ws.Range("O149").Value = 43
$ws.Range("P149").Value = 112
$ws.Range("Q149").Value = 135
$ws.Range("R149").Value = 141
$ws.Range("S149").Value = 224
$ws.Range("T149").Value = 41
$ws.Range("U149").Value = 443
$ws.Range("V149").Value = 214
$ws.Range("W149").Value = 510
$ws.Range("X149").Value = 636
$ws.Range("Y149").Value = 79
$ws.Range("Z149").Value = 497
$ws.Range("AA149").Value = 248
$ws.Range("AC149").Value = 2022
$ws.Range("AD149").Value = 43

$ws.Range("AE149").Formula = "=ROUND((P149-B149)/B149*100,2)"
$ws.Range("AF149").Formula = "=ROUND((Q149-C149)/C149*100,2)"
$ws.Range("AG149").Formula = "=ROUND((R149-D149)/D149*100,2)"
$ws.Range("AH149").Formula = "=ROUND((S149-E149)/E149*100,2)"
$ws.Range("AI149").Formula = "=ROUND((T149-F149)/F149*100,2)"
$ws.Range("AJ149").Formula = "=ROUND((U149-G149)/G149*100,2)"
$ws.Range("AK149").Formula = "=ROUND((V149-H149)/H149*100,2)"
$ws.Range("AL149").Formula = "=ROUND((W149-I149)/I149*100,2)"
$ws.Range("AM149").Formula = "=ROUND((X149-J149)/J149*100,2)"
$ws.Range("AN149").Formula = "=ROUND((Y149-K149)/K149*100,2)"
$ws.Range("AO149").Formula = "=ROUND((Z149-L149)/L149*100,2)"
$ws.Range("AP149").Formula = "=ROUND((AA149-M149)/M149*100,2)"

# --- Row 150 (new week 2022-44) --------------------------------------------
$ws.Range("N150").Value = 2022
$ws.Range("O150").Value = 44
$ws.Range("P150").Value = 108
$ws.Range("Q150").Value = 143
$ws.Range("R150").Value = 132
$ws.Range("S150").Value = 230
$ws.Range("T150").Value = 64
$ws.Range("U150").Value = 392
$ws.Range("V150").Value = 221
$ws.Range("W150").Value = 488
$ws.Range("X150").Value = 695
$ws.Range("Y150").Value = 89
$ws.Range("Z150").Value = 459
$ws.Range("AA150").Value = 222
$ws.Range("AC150").Value = 2022
$ws.Range("AD150").Value = 44

$ws.Range("AE150").Formula = "=ROUND((P150-B150)/B150*100,2)"
$ws.Range("AF150").Formula = "=ROUND((Q150-C150)/C150*100,2)"
$ws.Range("AG150").Formula = "=ROUND((R150-D150)/D150*100,2)"
$ws.Range("AH150").Formula = "=ROUND((S150-E150)/E150*100,2)"
$ws.Range("AI150").Formula = "=ROUND((T150-F150)/F150*100,2)"
$ws.Range("AJ150").Formula = "=ROUND((U150-G150)/G150*100,2)"
$ws.Range("AK150").Formula = "=ROUND((V150-H150)/H150*100,2)"
$ws.Range("AL150").Formula = "=ROUND((W150-I150)/I150*100,2)"
$ws.Range("AM150").Formula = "=ROUND((X150-J150)/J150*100,2)"
$ws.Range("AN150").Formula = "=ROUND((Y150-K150)/K150*100,2)"
$ws.Range("AO150").Formula = "=ROUND((Z150-L150)/L150*100,2)"
$ws.Range("AP150").Formula = "=ROUND((AA150-M150)/M150*100,2)"

# --- Restore the active selection/scroll position as recorded in the file --
$excel.ActiveWindow.ScrollRow = 82
$ws.Range("AH150").Select() | Out-Null
